# LanguageToolsTasks.xlsx - add "Note" column (F) with remarks for two tasks,
# matching the change where the LanguageWindow (only XML repository support)
# and the English-only keyboard feature got explanatory notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "Color the words based on the type in various languages" ---
# Progress % bumped from 0 to 70 and a note is added in column F.
# Set F6 before F1 so that "Works for English only" becomes shared-string
# index 18 (matching the target ordering).
$ws.Range("E6").Value = 70
$ws.Range("F6").Value = "Works for English only"

# --- Header row: new column F "Note" ---
# Give it the same "Good" (green) style used by the other header cells.
$ws.Range("F1").Value = "Note"
$ws.Range("F1").Style = "Good"

# --- Row 11: "Create and Edit Languague" ---
# Progress % bumped from 0 to 70 and a note is added in column F.
$ws.Range("E11").Value = 70
$ws.Range("F11").Value = "Only XMLRepository"

# Give new column F the same width style as the other data columns.
$ws.Columns("F").ColumnWidth = 25

# Update the selected cell to match the saved view state.
$ws.Range("F12").Select()
